$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91.0
$ws.Range("B3").Value = 46.0
$ws.Range("B4").Value = 72.0
$ws.Range("B6").Value = 478.0
